$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,20
$data[0,0] = "ECs"
$data[0,1] = "Cfh"
$data[0,2] = "Itgam"
$data[0,3] = "FAPs"
$data[0,4] = [double]"3"
$data[0,5] = [double]"1"
$data[0,6] = [double]"0.6292110000000001"
$data[0,7] = [double]"1.887633"
$data[0,8] = [double]"0.01078649253029594"
$data[0,9] = [double]"0.01078649253029594"
$data[0,10] = [double]"1"
$data[0,11] = [double]"0.3333333333333333"
$data[0,12] = [double]"0.142723"
$data[0,13] = [double]"0.428169"
$data[0,14] = [double]"0.0009642800942465787"
$data[0,15] = [double]"0.0009642800942465787"
$data[0,16] = [double]"0.08980288155300002"
$data[0,17] = [double]"0.808225933977"
$data[0,18] = [double]"1.040120003370378E-05"
$data[0,19] = [double]"1.040120003370378E-05"
$data[1,0] = "ECs"
$data[1,1] = "Cfh"
$data[1,2] = "Itgam"
$data[1,3] = "Inflammatory-Mac"
$data[1,4] = [double]"3"
$data[1,5] = [double]"1"
$data[1,6] = [double]"0.6292110000000001"
$data[1,7] = [double]"1.887633"
$data[1,8] = [double]"0.01078649253029594"
$data[1,9] = [double]"0.01078649253029594"
$data[1,10] = [double]"3"
$data[1,11] = [double]"1"
$data[1,12] = [double]"86.42780700000002"
$data[1,13] = [double]"259.283421"
$data[1,14] = [double]"0.5839326098770704"
$data[1,15] = [double]"0.5839326098770704"
$data[1,16] = [double]"54.38132687027701"
$data[1,17] = [double]"489.4319418324931"
$data[1,18] = [double]"0.006298584734635231"
$data[1,19] = [double]"0.006298584734635231"
$data[2,0] = "ECs"
$data[2,1] = "Cfh"
$data[2,2] = "Itgam"
$data[2,3] = "MuSCs"
$data[2,4] = [double]"3"
$data[2,5] = [double]"1"
$data[2,6] = [double]"0.6292110000000001"
$data[2,7] = [double]"1.887633"
$data[2,8] = [double]"0.01078649253029594"
$data[2,9] = [double]"0.01078649253029594"
$data[2,10] = [double]"1"
$data[2,11] = [double]"0.3333333333333333"
$data[2,12] = [double]"0.006361333333333333"
$data[2,13] = [double]"0.019084"
$data[2,14] = [double]"4.297910712499435E-05"
$data[2,15] = [double]"4.297910712499435E-05"
$data[2,16] = [double]"0.004002620908000001"
$data[2,17] = [double]"0.036023588172"
$data[2,18] = [double]"4.635938179625405E-07"
$data[2,19] = [double]"4.635938179625405E-07"
$data[3,0] = "ECs"
$data[3,1] = "Cfh"
$data[3,2] = "Itgam"
$data[3,3] = "Resolving-Mac"
$data[3,4] = [double]"3"
$data[3,5] = [double]"1"
$data[3,6] = [double]"0.6292110000000001"
$data[3,7] = [double]"1.887633"
$data[3,8] = [double]"0.01078649253029594"
$data[3,9] = [double]"0.01078649253029594"
$data[3,10] = [double]"3"
$data[3,11] = [double]"1"
$data[3,12] = [double]"61.43300833333333"
$data[3,13] = [double]"184.299025"
$data[3,14] = [double]"0.415060130921558"
$data[3,15] = [double]"0.415060130921558"
$data[3,16] = [double]"38.65432460642501"
$data[3,17] = [double]"347.888921457825"
$data[3,18] = [double]"0.004477043001809039"
$data[3,19] = [double]"0.004477043001809039"
$data[4,0] = "FAPs"
$data[4,1] = "Cfh"
$data[4,2] = "Itgam"
$data[4,3] = "FAPs"
$data[4,4] = [double]"3"
$data[4,5] = [double]"1"
$data[4,6] = [double]"24.05951033333333"
$data[4,7] = [double]"72.17853099999999"
$data[4,8] = [double]"0.4124494462002061"
$data[4,9] = [double]"0.4124494462002061"
$data[4,10] = [double]"1"
$data[4,11] = [double]"0.3333333333333333"
$data[4,12] = [double]"0.142723"
$data[4,13] = [double]"0.428169"
$data[4,14] = [double]"0.0009642800942465787"
$data[4,15] = [double]"0.0009642800942465787"
$data[4,16] = [double]"3.433845493304334"
$data[4,17] = [double]"30.904609439739"
$data[4,18] = [double]"0.000397716790853884"
$data[4,19] = [double]"0.0003977167908538839"
$data[5,0] = "FAPs"
$data[5,1] = "Cfh"
$data[5,2] = "Itgam"
$data[5,3] = "Inflammatory-Mac"
$data[5,4] = [double]"3"
$data[5,5] = [double]"1"
$data[5,6] = [double]"24.05951033333333"
$data[5,7] = [double]"72.17853099999999"
$data[5,8] = [double]"0.4124494462002061"
$data[5,9] = [double]"0.4124494462002061"
$data[5,10] = [double]"3"
$data[5,11] = [double]"1"
$data[5,12] = [double]"86.42780700000002"
$data[5,13] = [double]"259.283421"
$data[5,14] = [double]"0.5839326098770704"
$data[5,15] = [double]"0.5839326098770704"
$data[5,16] = [double]"2079.410715603839"
$data[5,17] = [double]"18714.69644043455"
$data[5,18] = [double]"0.2408426815620387"
$data[5,19] = [double]"0.2408426815620387"
$data[6,0] = "FAPs"
$data[6,1] = "Cfh"
$data[6,2] = "Itgam"
$data[6,3] = "MuSCs"
$data[6,4] = [double]"3"
$data[6,5] = [double]"1"
$data[6,6] = [double]"24.05951033333333"
$data[6,7] = [double]"72.17853099999999"
$data[6,8] = [double]"0.4124494462002061"
$data[6,9] = [double]"0.4124494462002061"
$data[6,10] = [double]"1"
$data[6,11] = [double]"0.3333333333333333"
$data[6,12] = [double]"0.006361333333333333"
$data[6,13] = [double]"0.019084"
$data[6,14] = [double]"4.297910712499435E-05"
$data[6,15] = [double]"4.297910712499435E-05"
$data[6,16] = [double]"0.1530505650671111"
$data[6,17] = [double]"1.377455085604"
$data[6,18] = [double]"1.772670893188325E-05"
$data[6,19] = [double]"1.772670893188325E-05"
$data[7,0] = "FAPs"
$data[7,1] = "Cfh"
$data[7,2] = "Itgam"
$data[7,3] = "Resolving-Mac"
$data[7,4] = [double]"3"
$data[7,5] = [double]"1"
$data[7,6] = [double]"24.05951033333333"
$data[7,7] = [double]"72.17853099999999"
$data[7,8] = [double]"0.4124494462002061"
$data[7,9] = [double]"0.4124494462002061"
$data[7,10] = [double]"3"
$data[7,11] = [double]"1"
$data[7,12] = [double]"61.43300833333333"
$data[7,13] = [double]"184.299025"
$data[7,14] = [double]"0.415060130921558"
$data[7,15] = [double]"0.415060130921558"
$data[7,16] = [double]"1478.048098803586"
$data[7,17] = [double]"13302.43288923227"
$data[7,18] = [double]"0.1711913211383816"
$data[7,19] = [double]"0.1711913211383816"
$data[8,0] = "Inflammatory-Mac"
$data[8,1] = "Cfh"
$data[8,2] = "Itgam"
$data[8,3] = "FAPs"
$data[8,4] = [double]"3"
$data[8,5] = [double]"1"
$data[8,6] = [double]"16.061552"
$data[8,7] = [double]"48.184656"
$data[8,8] = [double]"0.2753413571488099"
$data[8,9] = [double]"0.2753413571488098"
$data[8,10] = [double]"1"
$data[8,11] = [double]"0.3333333333333333"
$data[8,12] = [double]"0.142723"
$data[8,13] = [double]"0.428169"
$data[8,14] = [double]"0.0009642800942465787"
$data[8,15] = [double]"0.0009642800942465787"
$data[8,16] = [double]"2.292352886096001"
$data[8,17] = [double]"20.631175974864"
$data[8,18] = [double]"0.0002655061898214353"
$data[8,19] = [double]"0.0002655061898214352"
$data[9,0] = "Inflammatory-Mac"
$data[9,1] = "Cfh"
$data[9,2] = "Itgam"
$data[9,3] = "Inflammatory-Mac"
$data[9,4] = [double]"3"
$data[9,5] = [double]"1"
$data[9,6] = [double]"16.061552"
$data[9,7] = [double]"48.184656"
$data[9,8] = [double]"0.2753413571488099"
$data[9,9] = [double]"0.2753413571488098"
$data[9,10] = [double]"3"
$data[9,11] = [double]"1"
$data[9,12] = [double]"86.42780700000002"
$data[9,13] = [double]"259.283421"
$data[9,14] = [double]"0.5839326098770704"
$data[9,15] = [double]"0.5839326098770704"
$data[9,16] = [double]"1388.164716376464"
$data[9,17] = [double]"12493.48244738818"
$data[9,18] = [double]"0.1607807972869991"
$data[9,19] = [double]"0.1607807972869991"
$data[10,0] = "Inflammatory-Mac"
$data[10,1] = "Cfh"
$data[10,2] = "Itgam"
$data[10,3] = "MuSCs"
$data[10,4] = [double]"3"
$data[10,5] = [double]"1"
$data[10,6] = [double]"16.061552"
$data[10,7] = [double]"48.184656"
$data[10,8] = [double]"0.2753413571488099"
$data[10,9] = [double]"0.2753413571488098"
$data[10,10] = [double]"1"
$data[10,11] = [double]"0.3333333333333333"
$data[10,12] = [double]"0.006361333333333333"
$data[10,13] = [double]"0.019084"
$data[10,14] = [double]"4.297910712499435E-05"
$data[10,15] = [double]"4.297910712499435E-05"
$data[10,16] = [double]"0.1021728861226667"
$data[10,17] = [double]"0.9195559751040001"
$data[10,18] = [double]"1.183392568484003E-05"
$data[10,19] = [double]"1.183392568484003E-05"
$data[11,0] = "Inflammatory-Mac"
$data[11,1] = "Cfh"
$data[11,2] = "Itgam"
$data[11,3] = "Resolving-Mac"
$data[11,4] = [double]"3"
$data[11,5] = [double]"1"
$data[11,6] = [double]"16.061552"
$data[11,7] = [double]"48.184656"
$data[11,8] = [double]"0.2753413571488099"
$data[11,9] = [double]"0.2753413571488098"
$data[11,10] = [double]"3"
$data[11,11] = [double]"1"
$data[11,12] = [double]"61.43300833333333"
$data[11,13] = [double]"184.299025"
$data[11,14] = [double]"0.415060130921558"
$data[11,15] = [double]"0.415060130921558"
$data[11,16] = [double]"986.7094578622668"
$data[11,17] = [double]"8880.3851207604"
$data[11,18] = [double]"0.1142832197463045"
$data[11,19] = [double]"0.1142832197463045"
$data[12,0] = "MuSCs"
$data[12,1] = "Cfh"
$data[12,2] = "Itgam"
$data[12,3] = "FAPs"
$data[12,4] = [double]"3"
$data[12,5] = [double]"1"
$data[12,6] = [double]"1.855556"
$data[12,7] = [double]"5.566668"
$data[12,8] = [double]"0.03180958523221274"
$data[12,9] = [double]"0.03180958523221274"
$data[12,10] = [double]"1"
$data[12,11] = [double]"0.3333333333333333"
$data[12,12] = [double]"0.142723"
$data[12,13] = [double]"0.428169"
$data[12,14] = [double]"0.0009642800942465787"
$data[12,15] = [double]"0.0009642800942465787"
$data[12,16] = [double]"0.264830518988"
$data[12,17] = [double]"2.383474670892"
$data[12,18] = [double]"3.067334984566268E-05"
$data[12,19] = [double]"3.067334984566268E-05"
$data[13,0] = "MuSCs"
$data[13,1] = "Cfh"
$data[13,2] = "Itgam"
$data[13,3] = "Inflammatory-Mac"
$data[13,4] = [double]"3"
$data[13,5] = [double]"1"
$data[13,6] = [double]"1.855556"
$data[13,7] = [double]"5.566668"
$data[13,8] = [double]"0.03180958523221274"
$data[13,9] = [double]"0.03180958523221274"
$data[13,10] = [double]"3"
$data[13,11] = [double]"1"
$data[13,12] = [double]"86.42780700000002"
$data[13,13] = [double]"259.283421"
$data[13,14] = [double]"0.5839326098770704"
$data[13,15] = [double]"0.5839326098770704"
$data[13,16] = [double]"160.371635845692"
$data[13,17] = [double]"1443.344722611228"
$data[13,18] = [double]"0.0185746541237531"
$data[13,19] = [double]"0.0185746541237531"
$data[14,0] = "MuSCs"
$data[14,1] = "Cfh"
$data[14,2] = "Itgam"
$data[14,3] = "MuSCs"
$data[14,4] = [double]"3"
$data[14,5] = [double]"1"
$data[14,6] = [double]"1.855556"
$data[14,7] = [double]"5.566668"
$data[14,8] = [double]"0.03180958523221274"
$data[14,9] = [double]"0.03180958523221274"
$data[14,10] = [double]"1"
$data[14,11] = [double]"0.3333333333333333"
$data[14,12] = [double]"0.006361333333333333"
$data[14,13] = [double]"0.019084"
$data[14,14] = [double]"4.297910712499435E-05"
$data[14,15] = [double]"4.297910712499435E-05"
$data[14,16] = [double]"0.01180381023466667"
$data[14,17] = [double]"0.106234292112"
$data[14,18] = [double]"1.367147571296909E-06"
$data[14,19] = [double]"1.367147571296909E-06"
$data[15,0] = "MuSCs"
$data[15,1] = "Cfh"
$data[15,2] = "Itgam"
$data[15,3] = "Resolving-Mac"
$data[15,4] = [double]"3"
$data[15,5] = [double]"1"
$data[15,6] = [double]"1.855556"
$data[15,7] = [double]"5.566668"
$data[15,8] = [double]"0.03180958523221274"
$data[15,9] = [double]"0.03180958523221274"
$data[15,10] = [double]"3"
$data[15,11] = [double]"1"
$data[15,12] = [double]"61.43300833333333"
$data[15,13] = [double]"184.299025"
$data[15,14] = [double]"0.415060130921558"
$data[15,15] = [double]"0.415060130921558"
$data[15,16] = [double]"113.9923872109667"
$data[15,17] = [double]"1025.9314848987"
$data[15,18] = [double]"0.01320289061104267"
$data[15,19] = [double]"0.01320289061104267"
$data[16,0] = "Resolving-Mac"
$data[16,1] = "Cfh"
$data[16,2] = "Itgam"
$data[16,3] = "FAPs"
$data[16,4] = [double]"3"
$data[16,5] = [double]"1"
$data[16,6] = [double]"15.72740533333333"
$data[16,7] = [double]"47.182216"
$data[16,8] = [double]"0.2696131188884753"
$data[16,9] = [double]"0.2696131188884753"
$data[16,10] = [double]"1"
$data[16,11] = [double]"0.3333333333333333"
$data[16,12] = [double]"0.142723"
$data[16,13] = [double]"0.428169"
$data[16,14] = [double]"0.0009642800942465787"
$data[16,15] = [double]"0.0009642800942465787"
$data[16,16] = [double]"2.244662471389333"
$data[16,17] = [double]"20.201962242504"
$data[16,18] = [double]"0.000259982563691893"
$data[16,19] = [double]"0.000259982563691893"
$data[17,0] = "Resolving-Mac"
$data[17,1] = "Cfh"
$data[17,2] = "Itgam"
$data[17,3] = "Inflammatory-Mac"
$data[17,4] = [double]"3"
$data[17,5] = [double]"1"
$data[17,6] = [double]"15.72740533333333"
$data[17,7] = [double]"47.182216"
$data[17,8] = [double]"0.2696131188884753"
$data[17,9] = [double]"0.2696131188884753"
$data[17,10] = [double]"3"
$data[17,11] = [double]"1"
$data[17,12] = [double]"86.42780700000002"
$data[17,13] = [double]"259.283421"
$data[17,14] = [double]"0.5839326098770704"
$data[17,15] = [double]"0.5839326098770704"
$data[17,16] = [double]"1359.285152760104"
$data[17,17] = [double]"12233.56637484094"
$data[17,18] = [double]"0.1574358921696443"
$data[17,19] = [double]"0.1574358921696443"
$data[18,0] = "Resolving-Mac"
$data[18,1] = "Cfh"
$data[18,2] = "Itgam"
$data[18,3] = "MuSCs"
$data[18,4] = [double]"3"
$data[18,5] = [double]"1"
$data[18,6] = [double]"15.72740533333333"
$data[18,7] = [double]"47.182216"
$data[18,8] = [double]"0.2696131188884753"
$data[18,9] = [double]"0.2696131188884753"
$data[18,10] = [double]"1"
$data[18,11] = [double]"0.3333333333333333"
$data[18,12] = [double]"0.006361333333333333"
$data[18,13] = [double]"0.019084"
$data[18,14] = [double]"4.297910712499435E-05"
$data[18,15] = [double]"4.297910712499435E-05"
$data[18,16] = [double]"0.1000472677937778"
$data[18,17] = [double]"0.9004254101439999"
$data[18,18] = [double]"1.158773111901162E-05"
$data[18,19] = [double]"1.158773111901162E-05"
$data[19,0] = "Resolving-Mac"
$data[19,1] = "Cfh"
$data[19,2] = "Itgam"
$data[19,3] = "Resolving-Mac"
$data[19,4] = [double]"3"
$data[19,5] = [double]"1"
$data[19,6] = [double]"15.72740533333333"
$data[19,7] = [double]"47.182216"
$data[19,8] = [double]"0.2696131188884753"
$data[19,9] = [double]"0.2696131188884753"
$data[19,10] = [double]"3"
$data[19,11] = [double]"1"
$data[19,12] = [double]"61.43300833333333"
$data[19,13] = [double]"184.299025"
$data[19,14] = [double]"0.415060130921558"
$data[19,15] = [double]"0.415060130921558"
$data[19,16] = [double]"966.1818229043777"
$data[19,17] = [double]"8695.6364061394"
$data[19,18] = [double]"0.1119056564240201"
$data[19,19] = [double]"0.1119056564240201"

$ws.Range("A2:T21").Value = $data
Write-Output "done"
